# Update the "Table of Contents" title on the slide that contains it:
# translate the English title to Spanish (Ecuador) and tag the run's
# language accordingly, matching the author's edit.

$p = $ppt.ActivePresentation

$targetShape = $null

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "TABLE OF CONTENTS") {
                $targetShape = $shp
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange
$tr.Text = "TABLA DE CONTENIDOS"
$tr.LanguageID = "es-EC"
